$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row 5: "BOTON DE CERRAR SESISON Y QUE CUANDO INICIES APARESCA EL NOMBRE DEL USUARIO" ---
$ws.Rows(5).Insert()
# Pull formatting (thin-border style used by the data rows) from a plain row instead of
# re-typing style indices by hand, so the workbook's existing style table is reused as-is.
$ws.Range("B7:E7").Copy()
$ws.Range("B5:E5").PasteSpecial(-4122)
$ws.Range("B5").Value = "BOTON DE CERRAR SESISON Y QUE CUANDO INICIES APARESCA EL NOMBRE DEL USUARIO"

# --- Insert new row 9 (post first-insert numbering): "AJUSTAR EL LOGO DE LA PAGINA" ---
$ws.Rows(9).Insert()
$ws.Range("B7:E7").Copy()
$ws.Range("B9:E9").PasteSpecial(-4122)
$ws.Range("B9").Value = "AJUSTAR EL LOGO DE LA PAGINA"

$excel.CutCopyMode = 0

# Match the new selected cell recorded in the saved workbook
$ws.Range("B9").Select()
